# Actualización automática 2025-09-18 09:12:30
#
# Updates PORCELANATO sales figures for HIDALGO HIDALGO PEDRO GUSTAVO
# across the three report sheets, keeping the derived "total" rows and
# the dependent CUMPLIMIENTO MENSUAL (POR CUMPLIR / CUMPLIMIENTO) figures
# consistent with the new amounts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
#   Row 11 -> JARAMILLO CARVAJAL NICOLAS ESTEBAN, column M = PORCELANATO
#   Row 13 -> MEGAMAFERS S.A.,                    column M = PORCELANATO
# ---------------------------------------------------------------------
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("M11").Value = 5009.44
$wsVentasPorGrupo.Range("M13").Value = 4473.94

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
#   Row 11 / Row 13 hold the same two clients' monthly sales (column F).
#   Row 23 is the TOTAL row (sum of column F) and must be updated to
#   reflect the new figures.
# ---------------------------------------------------------------------
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F11").Value = 5925.64
$wsVentaMensual.Range("F13").Value = 4473.94
$wsVentaMensual.Range("F23").Value = 38324.91

# ---------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
#   Row 12 -> PORCELANATO group: VENTA (D), POR CUMPLIR (E = C - D) and
#             CUMPLIMIENTO (F = D / C) need to be refreshed.
#   Row 15 -> TOTAL row, same relationship.
#   Column E's width also reverts from 23 to 22.
# ---------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D12").Value = 34839.63
$wsCumplimiento.Range("E12").Value = 1984.0130921171
$wsCumplimiento.Range("F12").Value = 0.9461212165468271

$wsCumplimiento.Range("D15").Value = 38324.91
$wsCumplimiento.Range("E15").Value = 17099.83316613378
$wsCumplimiento.Range("F15").Value = 0.6914765465871874

# NOTE: Excel's Columns.ColumnWidth COM property is expressed in "characters"
# and is offset from the raw OOXML <col width="..."> units stored in the
# sheet XML (observed offset ~0.83 for this workbook's default font), so we
# compensate to land on the target raw width of 22.
$wsCumplimiento.Columns.Item(5).ColumnWidth = 21.17
